$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Narrow the last column ("AREA") from 1212 twips to 811 twips.
$t.Columns.Item(6).Width = 811 / 20.0

# Atlantic cod row: round the raw computed figures to 3 decimal places
# (Trange replaces T50 in the TEMPERATURE processing).
$t.Cell(2, 2).Range.Text = "1.122"
$t.Cell(2, 3).Range.Text = "0.000"
$t.Cell(2, 4).Range.Text = "2.000"
$t.Cell(2, 5).Range.Text = "0.561"
$t.Cell(2, 6).Range.Text = "0.561"

# European hake row: same rounding treatment.
$t.Cell(3, 2).Range.Text = "1.049"
$t.Cell(3, 3).Range.Text = "0.000"
$t.Cell(3, 4).Range.Text = "2.000"
$t.Cell(3, 5).Range.Text = "0.525"
$t.Cell(3, 6).Range.Text = "0.525"
